$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D36").Value = "What is Hard Negative Sample?"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/385"

$ws.Range("D42").Value = "[C#]셀프 프로세스 재실행"
$ws.Range("E42").Value = "https://kjk92.tistory.com/97"

$ws.Range("D51").Value = "[github] collaborator로 초대 받은 저장소 목록 확인하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/github-collaborator%EB%A1%9C-%EC%B4%88%EB%8C%80-%EB%B0%9B%EC%9D%80-%EC%A0%80%EC%9E%A5%EC%86%8C-%EB%AA%A9%EB%A1%9D-%ED%99%95%EC%9D%B8%ED%95%98%EA%B8%B0"
